$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds literal text in the source feed (e.g. "27.669.21",
# "311.60", "0.000008757", multi-dot thousands grouping, significant trailing zeros, etc.)
# rather than numbers. Force each target cell to Text format *before* writing so Excel
# does not auto-coerce the look-alike numerics and silently drop formatting (e.g.
# "311.60" -> 311.6). Style is restored to Normal afterwards to avoid leaving a stray
# Text-format style on cells that originally had the default style.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.669.21'
$ws.Range('E2').Value = '  -1.87%  '
$ws.Range('D2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.892.64'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('D3').Style = "Normal"

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.60'
$ws.Range('E5').Value = '  -1.01%  '
$ws.Range('D5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D6').Style = "Normal"

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4905'
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3795'
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('D8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07322'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('D9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9114'
$ws.Range('E10').Value = '  -4.16%  '
$ws.Range('D10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.58'
$ws.Range('E11').Value = '  -2.77%  '
$ws.Range('D11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '0.07651'
$ws.Range('E12').Value = '  -2.06%  '
$ws.Range('D12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.912.72'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D13').Style = "Normal"

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.461'
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.618'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.24'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008757'
$ws.Range('E18').Value = '  -2.06%  '
$ws.Range('D18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('D19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '27.708.61'
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('D20').Style = "Normal"

$ws.Range('E21').Value = '  -3.82%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.121'
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D22').Style = "Normal"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.128.27'
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('D23').Style = "Normal"

$ws.Range('E24').Value = '  -2.16%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.66'
$ws.Range('E25').Value = '  -2.18%  '
$ws.Range('D25').Style = "Normal"

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.862'
$ws.Range('E26').Value = '  -3.76%  '
$ws.Range('D26').Style = "Normal"

$ws.Range('E27').Value = '  -1.69%  '

$ws.Range('E28').Value = '  +3.37%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '115.16'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.871'
$ws.Range('E30').Value = '  -3.00%  '
$ws.Range('D30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08934'
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('D31').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.203'
$ws.Range('E32').Value = '  -3.76%  '
$ws.Range('D32').Style = "Normal"

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.222'
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('D33').Style = "Normal"

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7649'
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('D34').Style = "Normal"

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.620'
$ws.Range('E35').Value = '  -1.60%  '
$ws.Range('D35').Style = "Normal"

$ws.Range('D36').NumberFormat = "@"
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '0.02037'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('D36').Style = "Normal"

$ws.Range('D37').NumberFormat = "@"
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '2.548'
$ws.Range('E37').Value = '  -7.57%  '
$ws.Range('D37').Style = "Normal"

$ws.Range('E38').Value = '  -2.70%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05287'
$ws.Range('E39').Value = '  -2.04%  '
$ws.Range('D39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.5456'
$ws.Range('E40').Value = '  -3.49%  '
$ws.Range('D40').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.982'
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D41').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.886'
$ws.Range('E42').Value = '  -3.17%  '
$ws.Range('D42').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.533'
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D43').Style = "Normal"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '113.24'
$ws.Range('E44').Value = '  +6.42%  '
$ws.Range('D44').Style = "Normal"

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1520'
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('D45').Style = "Normal"

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.65'
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D46').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4783'
$ws.Range('E47').Value = '  -3.32%  '
$ws.Range('D47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D48').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.632'
$ws.Range('E49').Value = '  -3.23%  '
$ws.Range('D49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '67.46'
$ws.Range('E50').Value = '  -2.36%  '
$ws.Range('D50').Style = "Normal"

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06060'
$ws.Range('E51').Value = '  -1.19%  '
$ws.Range('D51').Style = "Normal"
